$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 34 (pushes existing rows 34..195 down to 35..196),
# then populate it with the new fruit record.
$ws.Rows.Item(34).Insert()
$ws.Range("A35").Copy($ws.Range("A34"))
$ws.Range("A34").Value = 34
$ws.Range("B34").Value = "grape"
$ws.Range("C34").Value = "Purple"
$ws.Range("D34").Value = "Tiny"
$ws.Range("E34").Value = 10.90219764355227

# Insert a second new data row at row 137 (in the now-shifted row numbering),
# pushing rows 137..196 down to 138..197, then populate it.
$ws.Rows.Item(137).Insert()
$ws.Range("A138").Copy($ws.Range("A137"))
$ws.Range("A137").Value = 138
$ws.Range("B137").Value = "grape"
$ws.Range("C137").Value = "Green"
$ws.Range("D137").Value = "Small"
$ws.Range("E137").Value = 10.90219764355227

Write-Host "edit complete"
